$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Rows 2-6 (2014/12 .. 2018/12 IFRS연결 columns): overwrite the
# reported figures with the corrected ones from the error-fix commit.
# J (당기순이익(비지배)) and O (자본총계(비지배)) are cleared for the
# rows that previously carried stray/incorrect minority-interest data.
# ------------------------------------------------------------------
    # Row 2
    $ws.Range("D2").Value = 1938
    $ws.Range("E2").Value = 164
    $ws.Range("F2").Value = 164
    $ws.Range("G2").Value = 201
    $ws.Range("H2").Value = 168
    $ws.Range("I2").Value = 168
    $ws.Range("J2").ClearContents()
    $ws.Range("K2").Value = 2416
    $ws.Range("L2").Value = 914
    $ws.Range("M2").Value = 1503
    $ws.Range("N2").Value = 1503
    $ws.Range("O2").ClearContents()
    $ws.Range("P2").Value = 64
    $ws.Range("Q2").Value = 194
    $ws.Range("R2").Value = -167
    $ws.Range("S2").Value = -21
    $ws.Range("T2").Value = 152
    $ws.Range("U2").Value = 42
    $ws.Range("V2").Value = 445
    $ws.Range("W2").Value = 8.45
    $ws.Range("X2").Value = 8.65
    $ws.Range("Y2").Value = 11.68
    $ws.Range("Z2").Value = 7.16
    $ws.Range("AA2").Value = 60.8
    $ws.Range("AB2").Value = 2383.36
    $ws.Range("AC2").Value = 1309
    $ws.Range("AD2").Value = 8.36
    $ws.Range("AE2").Value = 14310
    $ws.Range("AF2").Value = 0.77
    $ws.Range("AG2").Value = 200
    $ws.Range("AH2").Value = 1.83
    $ws.Range("AI2").Value = 12.53
    $ws.Range("AJ2").Value = 12800000
    # Row 3
    $ws.Range("D3").Value = 1895
    $ws.Range("E3").Value = 183
    $ws.Range("F3").Value = 183
    $ws.Range("G3").Value = 240
    $ws.Range("H3").Value = 193
    $ws.Range("I3").Value = 193
    $ws.Range("J3").ClearContents()
    $ws.Range("K3").Value = 2408
    $ws.Range("L3").Value = 747
    $ws.Range("M3").Value = 1661
    $ws.Range("N3").Value = 1661
    $ws.Range("O3").ClearContents()
    $ws.Range("P3").Value = 64
    $ws.Range("Q3").Value = 212
    $ws.Range("R3").Value = -30
    $ws.Range("S3").Value = -181
    $ws.Range("T3").Value = 29
    $ws.Range("U3").Value = 182
    $ws.Range("V3").Value = 285
    $ws.Range("W3").Value = 9.67
    $ws.Range("X3").Value = 10.18
    $ws.Range("Y3").Value = 12.2
    $ws.Range("Z3").Value = 8
    $ws.Range("AA3").Value = 44.99
    $ws.Range("AB3").Value = 2618.06
    $ws.Range("AC3").Value = 1508
    $ws.Range("AD3").Value = 9.22
    $ws.Range("AE3").Value = 15817
    $ws.Range("AF3").Value = 0.88
    $ws.Range("AG3").Value = 250
    $ws.Range("AH3").Value = 1.8
    $ws.Range("AI3").Value = 13.6
    $ws.Range("AJ3").Value = 12800000
    # Row 4
    $ws.Range("D4").Value = 1887
    $ws.Range("E4").Value = 167
    $ws.Range("F4").Value = 167
    $ws.Range("G4").Value = 218
    $ws.Range("H4").Value = 182
    $ws.Range("I4").Value = 182
    $ws.Range("J4").ClearContents()
    $ws.Range("K4").Value = 2524
    $ws.Range("L4").Value = 691
    $ws.Range("M4").Value = 1833
    $ws.Range("N4").Value = 1833
    $ws.Range("O4").ClearContents()
    $ws.Range("P4").Value = 64
    $ws.Range("Q4").Value = 182
    $ws.Range("R4").Value = -66
    $ws.Range("S4").Value = -107
    $ws.Range("T4").Value = 60
    $ws.Range("U4").Value = 121
    $ws.Range("V4").Value = 204
    $ws.Range("W4").Value = 8.83
    $ws.Range("X4").Value = 9.65
    $ws.Range("Y4").Value = 10.42
    $ws.Range("Z4").Value = 7.38
    $ws.Range("AA4").Value = 37.72
    $ws.Range("AB4").Value = 2861.9
    $ws.Range("AC4").Value = 1422
    $ws.Range("AD4").Value = 7.98
    $ws.Range("AE4").Value = 17454
    $ws.Range("AF4").Value = 0.65
    $ws.Range("AG4").Value = 400
    $ws.Range("AH4").Value = 3.52
    $ws.Range("AI4").Value = 23.07
    $ws.Range("AJ4").Value = 12800000
    # Row 5
    $ws.Range("D5").Value = 2012
    $ws.Range("E5").Value = 52
    $ws.Range("F5").Value = 52
    $ws.Range("G5").Value = 34
    $ws.Range("H5").Value = 43
    $ws.Range("I5").Value = 43
    $ws.Range("J5").ClearContents()
    $ws.Range("K5").Value = 2651
    $ws.Range("L5").Value = 832
    $ws.Range("M5").Value = 1818
    $ws.Range("N5").Value = 1818
    $ws.Range("O5").ClearContents()
    $ws.Range("P5").Value = 64
    $ws.Range("Q5").Value = 134
    $ws.Range("R5").Value = -268
    $ws.Range("S5").Value = 147
    $ws.Range("T5").Value = 193
    $ws.Range("U5").Value = -59
    $ws.Range("V5").Value = 393
    $ws.Range("W5").Value = 2.57
    $ws.Range("X5").Value = 2.14
    $ws.Range("Y5").Value = 2.36
    $ws.Range("Z5").Value = 1.66
    $ws.Range("AA5").Value = 45.78
    $ws.Range("AB5").Value = 2860.18
    $ws.Range("AC5").Value = 336
    $ws.Range("AD5").Value = 30.06
    $ws.Range("AE5").Value = 17317
    $ws.Range("AF5").Value = 0.58
    $ws.Range("AG5").Value = 250
    $ws.Range("AH5").Value = 2.48
    $ws.Range("AI5").Value = 61.04
    $ws.Range("AJ5").Value = 12800000
    # Row 6
    $ws.Range("D6").Value = 2095
    $ws.Range("E6").Value = 1
    $ws.Range("F6").Value = 1
    $ws.Range("G6").Value = -45
    $ws.Range("H6").Value = -39
    $ws.Range("I6").Value = -39
    $ws.Range("K6").Value = 3036
    $ws.Range("L6").Value = 1296
    $ws.Range("M6").Value = 1740
    $ws.Range("N6").Value = 1740
    $ws.Range("P6").Value = 64
    $ws.Range("Q6").Value = -50
    $ws.Range("R6").Value = -378
    $ws.Range("S6").Value = 420
    $ws.Range("T6").Value = 361
    $ws.Range("U6").Value = -411
    $ws.Range("V6").Value = 839
    $ws.Range("W6").Value = 0.06
    $ws.Range("X6").Value = -1.86
    $ws.Range("Y6").Value = -2.19
    $ws.Range("Z6").Value = -1.37
    $ws.Range("AA6").Value = 74.47
    $ws.Range("AB6").Value = 2759.91
    $ws.Range("AC6").Value = -305
    $ws.Range("AD6").Value = -21.42
    $ws.Range("AE6").Value = 16574
    $ws.Range("AF6").Value = 0.39
    $ws.Range("AG6").Value = 100
    $ws.Range("AH6").Value = 1.53
    $ws.Range("AI6").Value = -26.91
    $ws.Range("AJ6").Value = 12800000

# ------------------------------------------------------------------
# Rows 7-9 (2019/12(E) .. 2021/12(E) estimates): the erroneous
# forecast figures are removed entirely, leaving only the row number,
# "연간" label and the period label in columns A-C.
# ------------------------------------------------------------------
$ws.Range("D7:AJ9").ClearContents()
